{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer paragraphs\n// (plus the now-redundant blank paragraph that used to separate them from\n// the \"LOQ4003: ...\" requirement line), mirroring the site-rebuild diff\n// that dropped the scraped Jupiter/footer lines from the course page.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Identify the two footer paragraphs by their text (trimmed, in case of\n// stray whitespace differences picked up from the source document).\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nlet firstIdx = -1;\nlet lastIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (targetTexts.indexOf(t) !== -1) {\n    if (firstIdx === -1) {\n      firstIdx = i;\n    }\n    lastIdx = i;\n  }\n}\n\nif (firstIdx !== -1) {\n  // Also drop the blank paragraph that immediately precedes the block\n  // (the blank paragraph right after the \"\u00a9 2020 ...\" block is kept).\n  let startIdx = firstIdx;\n  if (firstIdx > 0 && items[firstIdx - 1].text.trim() === \"\") {\n    startIdx = firstIdx - 1;\n  }\n\n  const toDelete = [];\n  for (let i = startIdx; i <= lastIdx; i++) {\n    toDelete.push(items[i]);\n  }\n  // Delete in reverse so earlier deletions can't disturb later ones.\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer paragraphs\n# (plus the now-redundant blank paragraph that used to separate them from\n# the \"LOQ4003: ...\" requirement line), mirroring the site-rebuild diff\n# that dropped the scraped Jupiter/footer lines from the course page.\n\n$d = $word.ActiveDocument\n\n$target1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$target2 = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$firstIdx = -1\n$lastIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    $t = $t.TrimEnd([char]13, [char]7).Trim()\n    if ($t -eq $target1 -or $t -eq $target2) {\n        if ($firstIdx -eq -1) { $firstIdx = $i }\n        $lastIdx = $i\n    }\n}\n\nif ($firstIdx -ne -1) {\n    $startIdx = $firstIdx\n    if ($firstIdx -gt 1) {\n        $prev = $d.Paragraphs.Item($firstIdx - 1)\n        $prevText = $prev.Range.Text.TrimEnd([char]13, [char]7).Trim()\n        if ($prevText -eq \"\") {\n            $startIdx = $firstIdx - 1\n        }\n    }\n\n    $startPara = $d.Paragraphs.Item($startIdx)\n    $endPara = $d.Paragraphs.Item($lastIdx)\n    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $rng.Delete()\n}\n"}
